$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (REFERENCES -> CUSTOMERS); Excel auto-updates the
# sheet-scoped _xlnm._FilterDatabase defined name to follow the new name.
$ws.Name = "CUSTOMERS"

# --- Header-row style bookkeeping -----------------------------------------
# The header band (row 1 / row 2) keeps its exact visual borders, but the
# underlying style slots get rotated: B1:H1 <- I1 <- B2:H2 <- A2 <- B1:H1.
# Stage A2's current format in an unused scratch cell first so nothing is
# lost while the 4-way rotation is carried out with simple copy/paste of
# formats (keeps the same visual borders/fills throughout).
$ws.Range("A2").Copy()
$ws.Range("Z1").PasteSpecial(-4122)

$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("I1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("I1").PasteSpecial(-4122)

$ws.Range("Z1").Copy()
$ws.Range("B2:H2").PasteSpecial(-4122)

$ws.Range("Z1").Clear()

# Move the active selection to H17 to match the saved view state.
$ws.Range("H17").Select()
